$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "componentes"
$ws.Range("B1").Value = "cantidad"
$ws.Range("C1").Value = "fecha"
$ws.Range("D1").Value = "responsable"

# Data rows
$data = @(
    @("c1818", 25, "2024-01-23", "JORGE FUENTES"),
    @("c1818", 2,  "2024-01-23", "JORGE FUENTES"),
    @("c1818", 2,  "2024-01-23", "JORGE FUENTES"),
    @("c1818", 33, "2024-01-23", "JORGE FUENTES"),
    @("c1818", 3,  "2024-01-23", "JORGE FUENTES"),
    @("c1818", 3,  "2024-01-23", "JORGE FUENTES")
)

$ws.Range("C2:C7").NumberFormat = "@"

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $ws.Cells.Item($row, 4).Value = $item[3]
    $row++
}

# Drop the temporary text-number-format now that the literal date strings
# are locked in, so the data cells fall back to the default (unstyled) cell.
$ws.Range("C2:C7").ClearFormats()

# Header styling: bold, centered, top-aligned, thin box border
$headerRange = $ws.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Reset selection back to A1, like a freshly laid-out sheet
$ws.Range("A1").Select() | Out-Null
